$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds price strings that can look numeric (e.g. "513.24");
# force Text format on the whole column range first so COM stores the
# literal digits instead of silently parsing them into a Number, then
# drop the style back to Normal so no stray number-format style lingers
# on cells that did not have one before.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "56.776.81"
$ws.Range("E2").Value = "  -3.84%  "
$ws.Range("D3").Value = "2.534.65"
$ws.Range("E3").Value = "  -4.89%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").Value = "513.24"
$ws.Range("E5").Value = "  -2.18%  "
$ws.Range("D6").Value = "139.87"
$ws.Range("E6").Value = "  -3.21%  "
$ws.Range("E7").Value = "  +0.06%  "
$ws.Range("D8").Value = "0.556"
$ws.Range("E8").Value = "  -2.41%  "
$ws.Range("D9").Value = "6.53"
$ws.Range("E9").Value = "  -6.71%  "
$ws.Range("D10").Value = "0.0990"
$ws.Range("E10").Value = "  -3.93%  "
$ws.Range("D11").Value = "0.322"
$ws.Range("E11").Value = "  -3.88%  "
$ws.Range("E12").Value = "  -0.49%  "
$ws.Range("E13").Value = "  -4.68%  "
$ws.Range("D14").Value = "56.820.55"
$ws.Range("E14").Value = "  -3.76%  "
$ws.Range("D15").Value = "19.98"
$ws.Range("E15").Value = "  -5.22%  "
$ws.Range("E16").Value = "  -3.28%  "
$ws.Range("D17").Value = "2.543.57"
$ws.Range("E17").Value = "  -4.71%  "
$ws.Range("D18").Value = "330.44"
$ws.Range("E18").Value = "  -2.61%  "
$ws.Range("D19").Value = "4.26"
$ws.Range("E19").Value = "  -3.28%  "
$ws.Range("D20").Value = "10.06"
$ws.Range("E20").Value = "  -3.12%  "
$ws.Range("E21").Value = "  -4.44%  "
$ws.Range("E22").Value = "  +0.46%  "
$ws.Range("D23").Value = "64.98"
$ws.Range("E23").Value = "  +0.94%  "
$ws.Range("E24").Value = "  -0.66%  "
$ws.Range("D25").Value = "1.00"
$ws.Range("E25").Value = "  +0.18%  "
$ws.Range("D26").Value = "0.398"
$ws.Range("E26").Value = "  -4.86%  "
$ws.Range("D27").Value = "2.653.50"
$ws.Range("E27").Value = "  -4.64%  "
$ws.Range("E28").Value = "  -3.75%  "
$ws.Range("D29").Value = "0.0₃0746"
$ws.Range("E29").Value = "  -7.00%  "
$ws.Range("D30").Value = "0.998"
$ws.Range("E30").Value = "  -0.06%  "
$ws.Range("D31").Value = "6.23"
$ws.Range("E31").Value = "  -6.84%  "
$ws.Range("E32").Value = "  -3.14%  "
$ws.Range("D33").Value = "18.46"
$ws.Range("E33").Value = "  -2.21%  "
$ws.Range("D34").Value = "148.05"
$ws.Range("E34").Value = "  -1.69%  "
$ws.Range("D35").Value = "3.97"
$ws.Range("E35").Value = "  -4.47%  "
$ws.Range("E36").Value = "  -5.72%  "
$ws.Range("D37").Value = "0.842"
$ws.Range("E37").Value = "  -6.10%  "
$ws.Range("E38").Value = "  -3.63%  "
$ws.Range("D39").Value = "0.817"
$ws.Range("E39").Value = "  -6.70%  "
$ws.Range("D40").Value = "1.41"
$ws.Range("E40").Value = "  -3.08%  "
$ws.Range("E41").Value = "  +0.19%  "
$ws.Range("D42").Value = "3.46"
$ws.Range("E42").Value = "  -3.54%  "
$ws.Range("E43").Value = "  -2.02%  "
$ws.Range("E44").Value = "  -0.52%  "
$ws.Range("D45").Value = "0.576"
$ws.Range("E45").Value = "  -6.54%  "
$ws.Range("D46").Value = "260.43"
$ws.Range("E46").Value = "  -5.77%  "
$ws.Range("B47").Value = "EnergySwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D47").Value = "18.49"
$ws.Range("E47").Value = "  -7.12%  "
$ws.Range("B48").Value = "Hedera"
$ws.Range("C48").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D48").Value = "0.0515"
$ws.Range("E48").Value = "  -2.81%  "
$ws.Range("D49").Value = "1.964.61"
$ws.Range("E49").Value = "  -4.04%  "
$ws.Range("D50").Value = "0.0219"
$ws.Range("E50").Value = "  -4.31%  "
$ws.Range("D51").Value = "4.49"
$ws.Range("E51").Value = "  -4.73%  "

$ws.Range("D2:D51").Style = "Normal"
